# Updates the crypto price/volume table to the latest scraped values.
# Note: several "Price" cells (column D) are plain-looking numbers (e.g. "217.98")
# that Excel would otherwise auto-convert to a numeric type on assignment, changing
# the stored cell type from Text to Number (and losing exact formatting / producing
# float rounding). To keep those cells as Text - matching the source data - we
# assign them with a leading apostrophe (Excel's "store as text" marker) and then
# reset the cell Style back to "Normal" so no visible/formatting change is left
# behind on the cell itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.723.60'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.638.86'
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").Value = "'217.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("E6").Value = '  -0.83%  '

$ws.Range("E8").Value = '  -0.51%  '

$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("D10").Value = "'19.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.46%  '

$ws.Range("D11").Value = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").Value = '1.867.27'
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").Value = '1.631.83'
$ws.Range("E13").Value = '  -1.11%  '

$ws.Range("E14").Value = '  -1.27%  '

$ws.Range("E15").Value = '  -1.40%  '

$ws.Range("D16").Value = "'64.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.36%  '

$ws.Range("D17").Value = '26.696.87'
$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -2.10%  '

$ws.Range("D19").Value = "'211.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.11%  '

$ws.Range("E20").Value = '  -0.25%  '

$ws.Range("E21").Value = '  -0.60%  '

$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("E23").Value = '  -5.18%  '

$ws.Range("D24").Value = "'9.25"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'146.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("E27").Value = '  -1.62%  '

$ws.Range("E28").Value = '  -0.84%  '

$ws.Range("D29").Value = "'15.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.97%  '

$ws.Range("E30").Value = '  -3.19%  '

$ws.Range("D31").Value = "'1.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.47%  '

$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("D33").Value = "'2.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.62%  '

$ws.Range("D34").Value = '1.266.38'
$ws.Range("E34").Value = '  -1.15%  '

$ws.Range("E35").Value = '  -0.94%  '

$ws.Range("E36").Value = '  -0.80%  '

$ws.Range("E37").Value = '  -1.84%  '

$ws.Range("E38").Value = '  -1.71%  '

$ws.Range("D39").Value = "'0.803"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.91%  '

$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("D41").Value = "'0.802"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.45%  '

$ws.Range("D42").Value = "'2.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.17%  '

$ws.Range("D43").Value = '1.777.19'
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("D44").Value = "'5.26"
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("D46").Value = "'60.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.88%  '

$ws.Range("E47").Value = '  -1.70%  '

$ws.Range("E48").Value = '  +0.39%  '

$ws.Range("E49").Value = '  -3.11%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = "'0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.56%  '

Write-Host "Done applying changes"